# "Generate Report for Handback"
#
# The report workbook has three sheets:
#   - "Overview" : one row per source file, with a rolled-up
#                  "Latest HO Xliff Generate Date" column (G).
#   - "zh-cn"    : per-locale handoff/handback detail table (columns A-P).
#   - "de-de"    : same shape as "zh-cn", for the de-de locale.
#
# A new handback round just completed for the
# a33124e5-5ffc-40cc-b6c3-777194ec2b40.md source file (row 3 on every
# sheet) — the 7501bbdd... file (row 2) was already in sync and is left
# untouched. Refresh the handoff/handback timestamps for that row and
# roll the newest one up into the Overview sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# zh-cn: a33124e5 row (row 3) — new handoff / handback datetimes.
$wsZhCn.Range("H3").Value = "2016-09-07 07:08:21"   # Correspond Handoff Datetime
$wsZhCn.Range("K3").Value = "2016-09-07 07:08:56"   # Correspond Handback DateTime

# de-de: a33124e5 row (row 3) — new handoff / handback datetimes.
$wsDeDe.Range("H3").Value = "2016-09-07 07:08:27"   # Correspond Handoff Datetime
$wsDeDe.Range("K3").Value = "2016-09-07 07:09:12"   # Correspond Handback DateTime

# Overview: roll the newest "Latest HO Xliff Generate Date" up for the
# a33124e5 row (row 3).
$wsOverview.Range("G3").Value = "2016-09-07 07:08:27"
